$wb = $excel.ActiveWorkbook

# --- Sheet "Q4_19_20" (sheet1.xml) ---
$ws1 = $wb.Worksheets.Item("Q4_19_20")

# Update row 2 values
$ws1.Range("A2").Value = "Roads Places and Environment Group"
$ws1.Range("B2").Value = "SoT"

# Remove rows 3 through 7 (old data rows no longer present)
$ws1.Range("A3:C7").EntireRow.Delete()

# --- Sheet "Q4_18_19" (sheet2.xml) ---
$ws2 = $wb.Worksheets.Item("Q4_18_19")

# Add A2 value (was previously empty)
$ws2.Range("A2").Value = "None"

# Remove rows 3 through 6 (old data rows no longer present)
$ws2.Range("A3:C6").EntireRow.Delete()
